$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 5 de Octubre de 2020 a las 20:35"

# Helper to write a full data row (columns A-H)
function Set-Row($r, $country, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($r, 1).Value = $country
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
}

# Straight value refreshes (country/order unchanged)
Set-Row 4   "Estados Unidos" 7652925 16013 4864808 2573365 0 141 214752
Set-Row 5   "India"          6681953 59773 5657686 920667  0 886 103600
Set-Row 14  "Francia"        624274  5084  98680   493295  0 69  32299
Set-Row 34  "Marruecos"      134695  1423  113336  18990   0 39  2369
Set-Row 68  "Libano"         45657   1175  20243   25000   0 8   414
Set-Row 137 "Aruba"          4094    15    3612    451     0 1   31
Set-Row 143 "Sri Lanka"      3513    111   3259    241     0 0   13
Set-Row 145 "Mali"           3189    5     2482    576     0 0   131
Set-Row 180 "Curazao"        462     10    242     219     0 0   1
Set-Row 189 "Monaco"         223     1     193     28      0 0   2

# Peru / España swap places (España overtook Peru in total cases)
Set-Row 9   "España" 852838 2099 0      0     0 46 32225
Set-Row 10  "Peru"   828169 0    706223 89204 0 0  32742

# Costa Rica / Etiopia swap places (Etiopia overtook Costa Rica in total cases)
Set-Row 54  "Etiopia"    79437 618 34016 44191 0 8 1230
Set-Row 55  "Costa Rica" 79182 0   45007 33225 0 0 950
